$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Stamp the "About" sheet with the date this input file was generated/touched
# (2021-04-21, serial 44307) in a new column C, formatted as a short date.
$ws.Range("C1").Value = 44307
$ws.Range("C1").NumberFormat = "m/d/yyyy"
